# Scheduled market-data refresh for the Exodus_Profits workbook.
# Updates currentAveragePrice / Leve profit figures (columns H, I, J, K, L, M, N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets to match the
# latest Universalis price pull. A couple of rows gain/lose their profit
# (M) cell entirely where the recalculated NQ/HQ price relationship changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2096.4814
$ws.Range("J17").Value = 2171.0417
$ws.Range("L17").Value = 6513.125100000001
$ws.Range("N17").Value = -6849.125100000001

$ws.Range("H33").Value = 7550.826
$ws.Range("I33").Value = 7550.826
$ws.Range("K33").Value = 7550.826
$ws.Range("M33").Value = -7321.826

$ws.Range("H86").Value = 7491.552
$ws.Range("I86").Value = 7074.5454
$ws.Range("K86").Value = 7074.5454
$ws.Range("M86").Value = -5951.5454

$ws.Range("H89").Value = 7491.552
$ws.Range("I89").Value = 7074.5454
$ws.Range("K89").Value = 35372.727
$ws.Range("M89").Value = -29756.727

$ws.Range("H97").Value = 1745.2858
$ws.Range("J97").Value = 1745.2858
$ws.Range("L97").Value = 5235.857400000001
$ws.Range("N97").Value = -6227.857400000001

$ws.Range("H98").Value = 1178.7576
$ws.Range("I98").Value = 1178.7576
$ws.Range("K98").Value = 1178.7576
$ws.Range("M98").Value = 319.2424000000001

$ws.Range("H122").Value = 1178.7576
$ws.Range("I122").Value = 1178.7576
$ws.Range("K122").Value = 3536.2728
$ws.Range("M122").Value = -1086.2728

$ws.Range("H132").Value = 1442.9574
$ws.Range("I132").Value = 1240.4524
$ws.Range("K132").Value = 3721.357199999999
$ws.Range("M132").Value = -1191.357199999999

$ws.Range("H138").Value = 13515520
$ws.Range("J138").Value = 18520792
$ws.Range("L138").Value = 55562376
$ws.Range("N138").Value = -55572656

$ws.Range("H141").Value = 4622.5
$ws.Range("I141").Value = 4622.5
$ws.Range("K141").Value = 13867.5
$ws.Range("M141").Value = -8687.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 476
$ws.Range("I5").Value = 471.2
$ws.Range("K5").Value = 471.2
$ws.Range("M5").Value = -359.2

$ws.Range("H32").Value = 7582.1475
$ws.Range("I32").Value = 4427.3374
$ws.Range("J32").Value = 24407.8
$ws.Range("K32").Value = 4427.3374
$ws.Range("L32").Value = 24407.8
$ws.Range("M32").Value = -4140.3374
$ws.Range("N32").Value = -24981.8

$ws.Range("H45").Value = 53854.668
$ws.Range("I45").Value = 7033.8823
$ws.Range("K45").Value = 7033.8823
$ws.Range("M45").Value = -6656.8823

$ws.Range("H61").Value = 2491
$ws.Range("I61").Value = 1920.625
$ws.Range("K61").Value = 1920.625
$ws.Range("M61").Value = -1708.625

$ws.Range("H74").Value = 3390.25
$ws.Range("I74").Value = 2062.923
$ws.Range("J74").Value = 5855.2856
$ws.Range("K74").Value = 2062.923
$ws.Range("L74").Value = 5855.2856
$ws.Range("M74").Value = -1188.923
$ws.Range("N74").Value = -7603.2856

$ws.Range("H77").Value = 3390.25
$ws.Range("I77").Value = 2062.923
$ws.Range("J77").Value = 5855.2856
$ws.Range("K77").Value = 10314.615
$ws.Range("L77").Value = 29276.428
$ws.Range("M77").Value = -5946.614999999998
$ws.Range("N77").Value = -38012.428

$ws.Range("H132").Value = 2485.1714
$ws.Range("I132").Value = 2164.2727
$ws.Range("K132").Value = 6492.8181
$ws.Range("M132").Value = -3962.8181

$ws.Range("H136").Value = 2491
$ws.Range("I136").Value = 1920.625
$ws.Range("K136").Value = 5761.875
$ws.Range("M136").Value = -3211.875

$ws.Range("H139").Value = 92866.336
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 476
$ws.Range("I4").Value = 471.2
$ws.Range("K4").Value = 471.2
$ws.Range("M4").Value = -356.2

$ws.Range("H20").Value = 132075.23
$ws.Range("I20").Value = 171481.16
$ws.Range("J20").Value = 4006
$ws.Range("K20").Value = 171481.16
$ws.Range("L20").Value = 4006
$ws.Range("M20").Value = -171234.16
$ws.Range("N20").Value = -4500

$ws.Range("H134").Value = 2595.3057
$ws.Range("I134").Value = 2497.9707
$ws.Range("K134").Value = 7493.9121
$ws.Range("M134").Value = -4958.9121

$ws.Range("H138").Value = 99851
$ws.Range("J138").Value = 99851
$ws.Range("L138").Value = 99851
$ws.Range("N138").Value = -110131

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 145712.58
$ws.Range("I6").Value = 997.6
$ws.Range("K6").Value = 997.6
$ws.Range("M6").Value = -884.6

$ws.Range("H7").Value = 189.27586
$ws.Range("I7").Value = 119.15385
$ws.Range("J7").Value = 246.25
$ws.Range("K7").Value = 119.15385
$ws.Range("L7").Value = 246.25
$ws.Range("M7").Value = -6.153850000000006
$ws.Range("N7").Value = -472.25

$ws.Range("H22").Value = 429.07144
$ws.Range("I22").Value = 332
$ws.Range("J22").Value = 603.8
$ws.Range("K22").Value = 332
$ws.Range("L22").Value = 603.8
$ws.Range("M22").Value = 18
$ws.Range("N22").Value = -1303.8

$ws.Range("H58").Value = 1714.1
$ws.Range("I58").Value = 1315.1666
$ws.Range("K58").Value = 1315.1666
$ws.Range("M58").Value = -1112.1666

$ws.Range("H86").Value = 66669976
$ws.Range("I86").Value = 90912520
$ws.Range("K86").Value = 90912520
$ws.Range("M86").Value = -90911397

$ws.Range("H89").Value = 66669976
$ws.Range("I89").Value = 90912520
$ws.Range("K89").Value = 454562600
$ws.Range("M89").Value = -454556984

$ws.Range("H116").Value = 99960
$ws.Range("J116").Value = 99960
$ws.Range("L116").Value = 99960
$ws.Range("N116").Value = -109138

$ws.Range("H136").Value = 1714.1
$ws.Range("I136").Value = 1315.1666
$ws.Range("K136").Value = 3945.4998
$ws.Range("M136").Value = -1395.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1118.1428
$ws.Range("I36").Value = 956.75
$ws.Range("J36").Value = 1333.3334
$ws.Range("K36").Value = 2870.25
$ws.Range("L36").Value = 4000.0002
$ws.Range("M36").Value = -2701.25
$ws.Range("N36").Value = -4338.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3514
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 3119.6
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 3119.6
$ws.Range("M80").Value = -3502
$ws.Range("N80").Value = -5115.6

$ws.Range("H83").Value = 3514
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 3119.6
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 15598
$ws.Range("M83").Value = -17508
$ws.Range("N83").Value = -25582

$ws.Range("H97").Value = 802.4286
$ws.Range("I97").Value = 526
$ws.Range("J97").Value = 1171
$ws.Range("K97").Value = 526
$ws.Range("L97").Value = 1171
$ws.Range("M97").Value = -30
$ws.Range("N97").Value = -2163

$ws.Range("H111").Value = 67592.336
$ws.Range("J111").Value = 67592.336
$ws.Range("L111").Value = 67592.336
$ws.Range("M111").Value = -73726.336

$ws.Range("H126").Value = 3284.5588
$ws.Range("I126").Value = 2676.35
$ws.Range("J126").Value = 4153.4287
$ws.Range("K126").Value = 8029.049999999999
$ws.Range("L126").Value = 12460.2861
$ws.Range("M126").Value = -5559.049999999999
$ws.Range("N126").Value = -17400.2861

$ws.Range("H132").Value = 3673.4
$ws.Range("I132").Value = 3725.1
$ws.Range("J132").Value = 3570
$ws.Range("K132").Value = 11175.3
$ws.Range("L132").Value = 10710
$ws.Range("M132").Value = -8645.299999999999
$ws.Range("N132").Value = -15770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H40").Value = 6484911
$ws.Range("I40").Value = 3282.0833
$ws.Range("K40").Value = 3282.0833
$ws.Range("M40").Value = -3146.0833

$ws.Range("H46").Value = 1422.3636
$ws.Range("I46").Value = 1253
$ws.Range("J46").Value = 1718.75
$ws.Range("K46").Value = 1253
$ws.Range("L46").Value = 1718.75
$ws.Range("M46").Value = -1065
$ws.Range("N46").Value = -2094.75

$ws.Range("H97").Value = 7750
$ws.Range("J97").Value = 7750
$ws.Range("L97").Value = 7750
$ws.Range("N97").Value = -9732

$ws.Range("H122").Value = 10735284
$ws.Range("I122").Value = 29120.79
$ws.Range("J122").Value = 33337184
$ws.Range("K122").Value = 87362.37
$ws.Range("L122").Value = 100011552
$ws.Range("M122").Value = -84912.37
$ws.Range("N122").Value = -100016452

$ws.Range("H132").Value = 2578.96
$ws.Range("I132").Value = 2274.6428
$ws.Range("J132").Value = 2966.2727
$ws.Range("K132").Value = 6823.928400000001
$ws.Range("L132").Value = 8898.8181
$ws.Range("M132").Value = -4293.928400000001
$ws.Range("N132").Value = -13958.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2287.2114
$ws.Range("I122").Value = 1239.8462
$ws.Range("K122").Value = 3719.5386
$ws.Range("M122").Value = -1269.5386

$ws.Range("H132").Value = 1361477.2
$ws.Range("I132").Value = 2556.1052
$ws.Range("J132").Value = 3347593
$ws.Range("K132").Value = 7668.3156
$ws.Range("L132").Value = 10042779
$ws.Range("M132").Value = -5138.3156
$ws.Range("N132").Value = -10047839

$ws.Range("H138").Value = 93550.5
$ws.Range("J138").Value = 93550.5
$ws.Range("L138").Value = 93550.5
$ws.Range("N138").Value = -103830.5
